# Update the two-digit multiplication answers in the single table of the
# document. The table has 20 rows (5 "answer" rows containing the
# equations, interleaved with 4 blank rows each), 5 columns.
#
# Several old equation strings repeat (e.g. "72×14=1008" appears twice,
# in different cells, mapping to two different new values), so a global
# Find/Replace on $d.Content would be ambiguous. Instead we target each
# cell directly via the Tables collection and overwrite its Range.Text,
# which only touches that specific run and preserves its formatting
# (rFonts/sz) since Range.Text replaces the text node(s) in place.

$d   = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

# Row 1
$tbl.Cell(1, 1).Range.Text = "94×39=3666"   # was 13×74=962
$tbl.Cell(1, 2).Range.Text = "22×25=550"    # was 69×99=6831
$tbl.Cell(1, 3).Range.Text = "30×16=480"    # was 82×44=3608
$tbl.Cell(1, 4).Range.Text = "52×94=4888"   # was 97×71=6887
$tbl.Cell(1, 5).Range.Text = "82×96=7872"   # was 90×65=5850

# Row 5
$tbl.Cell(5, 1).Range.Text = "93×18=1674"   # was 83×52=4316
$tbl.Cell(5, 2).Range.Text = "12×55=660"    # was 72×14=1008
$tbl.Cell(5, 3).Range.Text = "63×66=4158"   # was 74×64=4736
$tbl.Cell(5, 4).Range.Text = "95×86=8170"   # was 63×69=4347
$tbl.Cell(5, 5).Range.Text = "98×16=1568"   # was 88×91=8008

# Row 10
$tbl.Cell(10, 1).Range.Text = "36×37=1332"  # was 74×59=4366
$tbl.Cell(10, 2).Range.Text = "84×93=7812"  # was 31×67=2077
$tbl.Cell(10, 3).Range.Text = "31×95=2945"  # was 66×85=5610
$tbl.Cell(10, 4).Range.Text = "70×73=5110"  # was 11×71=781
$tbl.Cell(10, 5).Range.Text = "15×71=1065"  # was 55×53=2915

# Row 15
$tbl.Cell(15, 1).Range.Text = "48×26=1248"  # was 35×97=3395
$tbl.Cell(15, 2).Range.Text = "56×47=2632"  # was 86×21=1806
$tbl.Cell(15, 3).Range.Text = "17×30=510"   # was 44×24=1056
$tbl.Cell(15, 4).Range.Text = "11×67=737"   # was 15×93=1395
$tbl.Cell(15, 5).Range.Text = "54×30=1620"  # was 98×77=7546

# Row 20
$tbl.Cell(20, 1).Range.Text = "71×20=1420"  # was 72×14=1008 (2nd occurrence)
$tbl.Cell(20, 2).Range.Text = "81×21=1701"  # was 76×77=5852
$tbl.Cell(20, 3).Range.Text = "58×54=3132"  # was 71×41=2911
$tbl.Cell(20, 4).Range.Text = "50×57=2850"  # was 95×92=8740
$tbl.Cell(20, 5).Range.Text = "83×65=5395"  # was 24×98=2352
